$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33, shifting existing rows 33-82 down to 34-83
$ws.Rows("33:33").Insert()

# Populate the new row 33 with the new data record
$ws.Range("A33").Value = 9
$ws.Range("B33").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C33").Value = "Metropolitana"
$ws.Range("D33").Value = 44797
$ws.Range("E33").Value = 13
$ws.Range("F33").Value = "Fruta"
$ws.Range("G33").Value = 100102
$ws.Range("H33").Value = "Cítricos"
$ws.Range("I33").Value = 100102006
$ws.Range("J33").Value = "Pomelo"
$ws.Range("K33").Value = "Start Ruby"
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 250
$ws.Range("N33").Value = 12000
$ws.Range("O33").Value = 12000
$ws.Range("P33").Value = 12000
$ws.Range("Q33").Value = "$/caja 14 kilos"
$ws.Range("R33").Value = "Región Metropolitana"
$ws.Range("S33").Value = 857
$ws.Range("T33").Value = 14
